# Add the "Dash Test App" and "NBA Prediction" project rows (rows 11 & 12)
# to the page_metadata worksheet, matching the source portfolio-website data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: dash_test_app -------------------------------------------------
$ws.Range("A11").Value = "dash_test_app"
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = "website"
$ws.Range("D11").Value = "Dash Test App"
$ws.Range("E11").Value = "https://github.com/cdpeters/dash-test-app"
$ws.Range("F11").Value = "Python, CSS"
$ws.Range("G11").Value = "pandas, dash, ibis-framework, sqlite, tailwindcss"
$ws.Range("H11").Value = "database, app"

# --- Row 12: nba_prediction -------------------------------------------------
$ws.Range("A12").Value = "nba_prediction"
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = "website"
$ws.Range("D12").Value = "NBA Prediction"
$ws.Range("E12").Value = "https://github.com/pascalegeday/NBA_Prediction_Analysis"
$ws.Range("F12").Value = "Python, JavaScript, HTML, CSS"
$ws.Range("G12").Value = "pandas, splinter, bs4 (beautiful soup 4), sqlalchemy, bootstrap, AWS, postgresql, flask "
$ws.Range("H12").Value = "database, web scraping, app, logistic regression, PCA"

# --- Formatting: match the other data rows (font id 1 = size 13 Calibri) ---
# (column E is excluded here - it keeps the separate "Hyperlink" style/font
# applied further down, same as E3:E10 above it)
$ws.Range("A11:D12").Font.Size = 13
$ws.Range("F11:H12").Font.Size = 13
$ws.Range("B11:B12").NumberFormat = "0"

# Row 11 specific alignment tweaks (new cellXfs 8/9/10 in the target file)
$ws.Range("A11").VerticalAlignment = -4108   # xlVAlignCenter
$ws.Range("D11").VerticalAlignment = -4108
$ws.Range("F11").VerticalAlignment = -4108
$ws.Range("G11").VerticalAlignment = -4108

$ws.Range("C11").VerticalAlignment = -4108
$ws.Range("C11").HorizontalAlignment = -4131  # xlHAlignLeft

$ws.Range("H11").VerticalAlignment = -4108
$ws.Range("H11").WrapText = $true

# --- Hyperlinks (added last) -------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E11"), "https://github.com/cdpeters/dash-test-app") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E12"), "https://github.com/pascalegeday/NBA_Prediction_Analysis") | Out-Null

# Adding a hyperlink forces a fresh (font-duplicated) style onto the cell;
# re-apply the already-existing "Hyperlink" cell style (same one used by
# E3:E10) so no redundant cellXfs entry gets appended.
$ws.Range("E11").Style = $ws.Range("E10").Style
$ws.Range("E12").Style = $ws.Range("E10").Style

# --- Selection matches the saved workbook state ------------------------------
$ws.Range("B12").Select()
